$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Reuse the exact same header style as the existing columns (H1, etc.)
# by copying formats from an existing header cell, rather than building a
# new composite style.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("I2").Value = 9
$ws.Range("J2").Value = 9
